# "Generate Report for Handoff"
# A new localization handoff cycle was generated for the source file that used
# to be named 5c93d2ed-78d5-4851-8388-4277466e853a.md; it has been renamed to
# d75331d2-192a-4f77-886a-577d3632db46.md, fresh xliff hand-off files were
# produced (new content hash 6001d09275f1457c88422eaec61b055a06a10ab6), and
# since the new cycle has not been handed *back* yet, the "target"/"handback"
# info for each locale is reset to "not yet available".

$wb = $excel.ActiveWorkbook

$oldBase = "5c93d2ed-78d5-4851-8388-4277466e853a"
$newBase = "d75331d2-192a-4f77-886a-577d3632db46"
$newHash = "6001d09275f1457c88422eaec61b055a06a10ab6"

# Same external target on every sheet's "source file" hyperlink (unchanged by
# this edit - the underlying link keeps pointing at the commit where the file
# still had its old name; only the visible/display text moves to the new
# file name).
$sourceUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/5717398e99f1c9bc6429217c27e0dfb156046fec/e2e/$oldBase.md"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newBase.md"
$wsOverview.Range("G2").Value = "2016-09-04 15:03:34"

# Re-point the hyperlink on B2 at the new display text (re-add, since
# mutating TextToDisplay in place on an existing Hyperlink object does not
# update the underlying XML in this runtime).
$wsOverview.Range("B2").Value = "e2e\$newBase.md"
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), $sourceUrl, [Type]::Missing, [Type]::Missing, "e2e\$newBase.md") | Out-Null

# ---------------------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
$locales = @(
    @{ Sheet = "zh-cn"; Suffix = "zh-cn"; HandoffDate = "2016-09-04 15:03:30" },
    @{ Sheet = "de-de"; Suffix = "de-de"; HandoffDate = "2016-09-04 15:03:34" }
)

foreach ($loc in $locales) {
    $ws = $wb.Worksheets.Item($loc.Sheet)

    $ws.Range("A2").Value = "$newBase.md"
    $ws.Range("G2").Value = "$newBase.$newHash.$($loc.Suffix).xlf"
    $ws.Range("H2").Value = $loc.HandoffDate

    # No handback has happened yet for the new cycle - target/handback file
    # columns go blank, and the handback timestamp resets to the zero date.
    $ws.Range("I2").Value = ""
    $ws.Range("I2").Style = "Normal"
    $ws.Range("J2").Value = ""
    $ws.Range("K2").Value = "0001-01-01 00:00:00"

    # Drop every hyperlink on the sheet (A2 + I2) and re-add only the one on
    # A2, now pointing at the new file's display text.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $sourceUrl, [Type]::Missing, [Type]::Missing, "$newBase.md") | Out-Null

    # Columns I/J shrink now that they hold short/empty values instead of
    # long file names.
    $ws.Columns.Item(9).ColumnWidth = 18.6506053379604
    $ws.Columns.Item(10).ColumnWidth = 21.7054770333426
}
